$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "TestDataMappingSheet_SD"

# Change the selected cell to A18
$ws.Range("A18").Select()
